# Mark the "Done" column (C) for the tasks that are now complete:
# Web App rows (23-26), part of CI/CD (29), Deploy (31-32) use "Y",
# while CI/CD row 29 and 30 use lowercase "y" (matching the existing
# pattern used for rows 17-20), and Submission rows (35-36) use "Y".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C23").Value = "Y"
$ws.Range("C24").Value = "Y"
$ws.Range("C25").Value = "Y"
$ws.Range("C26").Value = "Y"

$ws.Range("C29").Value = "y"
$ws.Range("C30").Value = "y"

$ws.Range("C31").Value = "Y"
$ws.Range("C32").Value = "Y"

$ws.Range("C35").Value = "Y"
$ws.Range("C36").Value = "Y"

# Update the view's visible top-left cell and active selection to
# reflect scrolling down to the bottom of the now-completed checklist
# (was topLeftCell A17 / selection C28, now A22 / C36).
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 22
$ws.Range("C36").Select()
